$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new timestamp row (2026/02/20 金 5:00) was inserted right before the
# existing "2026/12/29" block, shifting every row below it down by one
# (old row 834 -> new row 835, ..., old row 875 -> new row 876).
$ws.Rows(834).Insert()

# Leading apostrophe forces the date-looking string to stay plain text
# instead of Excel auto-coercing it into a date serial number; ClearFormats
# afterwards drops the transient "quote prefix" cell style so the new row
# ends up with the same (default/no) style as its neighbours.
$ws.Range("A834").Value = "'2026/02/20"
$ws.Range("A834").ClearFormats()

$ws.Range("B834").Value = "金"
$ws.Range("C834").Value = 5
$ws.Range("D834").Value = 201
